# Update values in column A per the diff (Update Name of Algo - KNN re-run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12 = -21.626
    32 = -21.648
    36 = -20.339
    38 = -19.741
    46 = -21.882
    54 = -22.209
    55 = -22.164
    67 = -21.495
    69 = -21.721
    72 = -21.445
    91 = -21.522
    99 = -20.43
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value = $updates[$row]
}
